# TLWP-990 - fix provider name on spreadsheet
# Remove the merged "Showing providers with courses in these skill areas: "
# header row (row 1) together with the blank spacer row (row 2) that sat
# above the real column-header row, so the real header row becomes row 1
# and a blank line remains after it (old row 2 gap is preserved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old filter-header row and the blank row beneath it; this
# shifts the column-header row (previously row 3) up to row 1, leaving
# row 2 blank again.
$ws.Range("A1:A2").EntireRow.Delete()

# Match the author's resulting selection (active cell on the blank row
# just below the new header row).
$ws.Range("A2").Select()
